$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2:F2) - date shifted, B/C/F cleared, D/E updated
$ws.Range("A2").Value = 45919
$ws.Range("D2").Value = 5592
$ws.Range("E2").Value = 6515.851061
$ws.Range("B2:C2").ClearContents()
$ws.Range("F2").ClearContents()

# Row 3 (A3:F3)
$ws.Range("A3").Value = 45920
$ws.Range("D3").Value = 2952
$ws.Range("E3").Value = 2320.793555
$ws.Range("B3:C3").ClearContents()
$ws.Range("F3").ClearContents()

# Row 4 (A4:F4)
$ws.Range("A4").Value = 45921
$ws.Range("D4").Value = 2952
$ws.Range("E4").Value = 2202.911293
$ws.Range("B4:C4").ClearContents()
$ws.Range("F4").ClearContents()

# Row 5 (A5:F5)
$ws.Range("A5").Value = 45922
$ws.Range("D5").Value = 2952
$ws.Range("E5").Value = 6527.061986
$ws.Range("B5:C5").ClearContents()
$ws.Range("F5").ClearContents()

# Re-apply a no-op font setting over the cleared cells so Excel keeps them as
# bare empty <c r=".."/> elements without allocating a new style index.
$ws.Range("B2:C5").Font.Bold = $false
$ws.Range("F2:F5").Font.Bold = $false

# Row 6 (A6:F6)
$ws.Range("A6").Value = 45923
$ws.Range("B6").Value = 5796.01973027908
$ws.Range("C6").Value = 5182.45604544916
$ws.Range("D6").Value = 2952
$ws.Range("E6").Value = 6523.437267
$ws.Range("F6").Value = 123.24473259042

# Row 7 (A7:F7)
$ws.Range("A7").Value = 45924
$ws.Range("B7").Value = 5979.99309220113
$ws.Range("C7").Value = 5456.5001150602
$ws.Range("D7").Value = 2952
$ws.Range("E7").Value = 6709.394895
$ws.Range("F7").Value = 134.745913244128

# Row 8 (A8:F8)
$ws.Range("A8").Value = 45925
$ws.Range("B8").Value = 5979.99309220113
$ws.Range("C8").Value = 5764.12867090187
$ws.Range("D8").Value = 2952
$ws.Range("E8").Value = 6709.394895
$ws.Range("F8").Value = 147.563769737531

# Row 9 (A9:F9)
$ws.Range("A9").Value = 45926
$ws.Range("B9").Value = 5979.99309220113
$ws.Range("C9").Value = 5059.27485981321
$ws.Range("D9").Value = 2952
$ws.Range("E9").Value = 6709.394895
$ws.Range("F9").Value = 118.19486094217

# Row 10 (A10:F10)
$ws.Range("A10").Value = 45927
$ws.Range("B10").Value = 1198.04914556408
$ws.Range("C10").Value = 2517.94701991763
$ws.Range("D10").Value = 2952
$ws.Range("E10").Value = 2322.416874
$ws.Range("F10").Value = 28.7631145147311

# Row 11 (A11:F11)
$ws.Range("A11").Value = 45928
$ws.Range("B11").Value = 1070.29087265608
$ws.Range("C11").Value = 2683.88342316486
$ws.Range("D11").Value = 2952
$ws.Range("E11").Value = 2183.81959
$ws.Range("F11").Value = 35.2255058545327

# Row 12 (A12:F12)
$ws.Range("A12").Value = 45929
$ws.Range("B12").Value = 6235.87706540818
$ws.Range("C12").Value = 5185.87877444631
$ws.Range("D12").Value = 2952
$ws.Range("E12").Value = 7030.796743
$ws.Range("F12").Value = 126.199935501589

# Row 13 (A13:F13)
$ws.Range("A13").Value = 45930
$ws.Range("B13").Value = 6235.87706540818
$ws.Range("C13").Value = 5356.54355997924
$ws.Range("D13").Value = 2952
$ws.Range("E13").Value = 7030.796743
$ws.Range("F13").Value = 133.310968232127

# Row 14 (A14:F14)
$ws.Range("A14").Value = 45931
$ws.Range("B14").Value = 5123.60662889003
$ws.Range("C14").Value = 4575.99821602637
$ws.Range("D14").Value = 3692
$ws.Range("E14").Value = 6097.869181
$ws.Range("F14").Value = 77.4275320056812

# Row 15 (A15:F15)
$ws.Range("A15").Value = 45932
$ws.Range("B15").Value = 5123.60662889003
$ws.Range("C15").Value = 4559.73383494078
$ws.Range("D15").Value = 3692
$ws.Range("E15").Value = 6097.869181
$ws.Range("F15").Value = 76.7498494604482

Write-Host "edit applied"
